# Applies the scheduled Leve-profit recalculation update to all 8 job sheets.
# Values below come from the upstream scraper refresh (commit: "chore: update Sheets via scheduled runner").
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 5
$ws.Range("H5").Value = 836.5714
$ws.Range("I5").Value = 131.75
$ws.Range("J5").Value = 1776.3334
$ws.Range("K5").Value = 131.75
$ws.Range("L5").Value = 1776.3334
$ws.Range("M5").Value = -16.75
$ws.Range("N5").Value = -2006.3334

# Row 17
$ws.Range("H17").Value = 1337.3043
$ws.Range("J17").Value = 1337.3043
$ws.Range("L17").Value = 4011.9129
$ws.Range("N17").Value = -4347.9129

# Row 69
$ws.Range("H69").Value = 12749
$ws.Range("I69").Value = 12331.667
$ws.Range("J69").Value = 13166.333
$ws.Range("K69").Value = 36995.001
$ws.Range("L69").Value = 39498.999
$ws.Range("M69").Value = -36121.001
$ws.Range("N69").Value = -41246.999

# Row 72
$ws.Range("H72").Value = 12749
$ws.Range("I72").Value = 12331.667
$ws.Range("J72").Value = 13166.333
$ws.Range("K72").Value = 110985.003
$ws.Range("L72").Value = 118496.997
$ws.Range("M72").Value = -106617.003
$ws.Range("N72").Value = -127232.997

# Row 129
$ws.Range("H129").Value = 2477.524
$ws.Range("I129").Value = 249.5
$ws.Range("J129").Value = 2712.0527
$ws.Range("K129").Value = 748.5
$ws.Range("L129").Value = 8136.158100000001
$ws.Range("M129").Value = 4251.5
$ws.Range("N129").Value = -18136.1581

# Row 137
$ws.Range("H137").Value = 1699.25
$ws.Range("J137").Value = 1999.5
$ws.Range("L137").Value = 5998.5
$ws.Range("N137").Value = -11098.5

# Row 138
$ws.Range("H138").Value = 4287.303
$ws.Range("J138").Value = 5099.8696
$ws.Range("L138").Value = 15299.6088
$ws.Range("N138").Value = -25579.6088

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 1986.25
$ws.Range("I61").Value = 1986.25
$ws.Range("K61").Value = 1986.25
$ws.Range("M61").Value = -1774.25

# Row 136
$ws.Range("H136").Value = 1986.25
$ws.Range("I136").Value = 1986.25
$ws.Range("K136").Value = 5958.75
$ws.Range("M136").Value = -3408.75

$ws = $wb.Worksheets.Item("BSM")
# Row 5
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("N5").ClearContents()

# Row 64
$ws.Range("H64").Value = 823.375
$ws.Range("J64").Value = 749.4
$ws.Range("L64").Value = 749.4
$ws.Range("N64").Value = -1199.4

# Row 67
$ws.Range("H67").Value = 823.375
$ws.Range("J67").Value = 749.4
$ws.Range("L67").Value = 749.4
$ws.Range("N67").Value = -2309.4

# Row 100
$ws.Range("H100").Value = 150000
$ws.Range("J100").Value = 150000
$ws.Range("L100").Value = 150000
$ws.Range("N100").Value = -152164

# Row 140
$ws.Range("H140").Value = 100000
$ws.Range("J140").Value = 100000
$ws.Range("L140").Value = 100000
$ws.Range("N140").Value = -110360

$ws = $wb.Worksheets.Item("CRP")
# Row 5
$ws.Range("H5").Value = 1181.75
$ws.Range("I5").Value = 1481
$ws.Range("J5").Value = 683
$ws.Range("K5").Value = 1481
$ws.Range("L5").Value = 683
$ws.Range("M5").Value = -1369
$ws.Range("N5").Value = -907

# Row 12
$ws.Range("H12").Value = 4992.5
$ws.Range("I12").Value = 4985
$ws.Range("J12").Value = 5000
$ws.Range("K12").Value = 4985
$ws.Range("L12").Value = 5000
$ws.Range("M12").Value = -4815
$ws.Range("N12").Value = -5340

# Row 22
$ws.Range("H22").Value = 332.42856
$ws.Range("I22").Value = 315.6
$ws.Range("K22").Value = 315.6
$ws.Range("M22").Value = 34.39999999999998

# Row 132
$ws.Range("H132").Value = 2592.6
$ws.Range("I132").Value = 2706.3572
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 8119.071599999999
$ws.Range("L132").Value = 3000
$ws.Range("M132").Value = -5589.071599999999
$ws.Range("N132").Value = -8060

# Row 141
$ws.Range("H141").Value = 152884.11
$ws.Range("I141").Value = 40000
$ws.Range("J141").Value = 166994.62
$ws.Range("K141").Value = 40000
$ws.Range("L141").Value = 166994.62
$ws.Range("M141").Value = -34820
$ws.Range("N141").Value = -177354.62

$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 2003
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 2003
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 6009
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -7631

# Row 71
$ws.Range("H71").Value = 2003
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 2003
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 18027
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -26139

# Row 103
$ws.Range("H103").Value = 946.75
$ws.Range("I103").Value = 45.666668
$ws.Range("J103").Value = 1487.4
$ws.Range("K103").Value = 137.000004
$ws.Range("L103").Value = 4462.200000000001
$ws.Range("M103").Value = 741.999996
$ws.Range("N103").Value = -6220.200000000001

# Row 129
$ws.Range("H129").Value = 3603
$ws.Range("I129").Value = 1964.8
$ws.Range("J129").Value = 6333.3335
$ws.Range("K129").Value = 5894.4
$ws.Range("L129").Value = 19000.0005
$ws.Range("M129").Value = -894.3999999999996
$ws.Range("N129").Value = -29000.0005

# Row 132
$ws.Range("H132").Value = 6315.6665
$ws.Range("J132").Value = 7632.3335
$ws.Range("L132").Value = 68691.0015
$ws.Range("N132").Value = -73751.0015

# Row 134
$ws.Range("H134").Value = 1850
$ws.Range("I134").Value = 1850
$ws.Range("K134").Value = 5550
$ws.Range("M134").Value = -480

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 84
$ws.Range("I2").Value = 96.59090999999999
$ws.Range("J2").Value = 14.75
$ws.Range("K2").Value = 96.59090999999999
$ws.Range("L2").Value = 14.75
$ws.Range("M2").Value = 16.40909000000001
$ws.Range("N2").Value = -240.75

# Row 10
$ws.Range("H10").Value = 14166.667
$ws.Range("J10").Value = 15250
$ws.Range("L10").Value = 15250
$ws.Range("N10").Value = -15588

# Row 21
$ws.Range("H21").Value = 8000
$ws.Range("J21").Value = 8000
$ws.Range("L21").Value = 8000
$ws.Range("N21").Value = -8346

# Row 30
$ws.Range("H30").Value = 8000
$ws.Range("J30").Value = 8000
$ws.Range("L30").Value = 8000
$ws.Range("N30").Value = -8210

# Row 132
$ws.Range("H132").Value = 4998
$ws.Range("I132").Value = 4998
$ws.Range("K132").Value = 14994
$ws.Range("M132").Value = -12464

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 512.3333
$ws.Range("I16").Value = 512.3333
$ws.Range("K16").Value = 512.3333
$ws.Range("M16").Value = -342.3333

# Row 22
$ws.Range("H22").Value = 1938.5834
$ws.Range("I22").Value = 1530.6666
$ws.Range("K22").Value = 1530.6666
$ws.Range("M22").Value = -1235.6666

# Row 27
$ws.Range("H27").Value = 1938.5834
$ws.Range("I27").Value = 1530.6666
$ws.Range("K27").Value = 1530.6666
$ws.Range("M27").Value = -1423.6666

# Row 93
$ws.Range("H93").Value = 2751.818
$ws.Range("I93").Value = 2561.111
$ws.Range("J93").Value = 3610
$ws.Range("K93").Value = 2561.111
$ws.Range("L93").Value = 3610
$ws.Range("M93").Value = -1313.111
$ws.Range("N93").Value = -6106

# Row 136
$ws.Range("H136").Value = 2913.1428
$ws.Range("I136").Value = 1466.6666
$ws.Range("K136").Value = 4399.9998
$ws.Range("M136").Value = -1849.9998

$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 4924.2163
$ws.Range("I126").Value = 4679.696
$ws.Range("J126").Value = 5325.9287
$ws.Range("K126").Value = 14039.088
$ws.Range("L126").Value = 15977.7861
$ws.Range("M126").Value = -11569.088
$ws.Range("N126").Value = -20917.7861
